$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(4, 1).Value = -20.84359999999998
$ws.Cells.Item(6, 1).Value = -22.79040000000001
$ws.Cells.Item(7, 1).Value = -21.9261
$ws.Cells.Item(7, 2).Value = 4.660599999999999
$ws.Cells.Item(12, 2).Value = 4.738599999999999
$ws.Cells.Item(13, 5).Value = 16.73
$ws.Cells.Item(14, 5).Value = 16.9919
$ws.Cells.Item(15, 2).Value = 5.294699999999997
$ws.Cells.Item(16, 1).Value = -21.93480000000001
$ws.Cells.Item(16, 5).Value = 16.74079999999999
$ws.Cells.Item(19, 5).Value = 16.35959999999999
$ws.Cells.Item(20, 1).Value = -23.01780000000001
$ws.Cells.Item(20, 2).Value = 4.761899999999998
$ws.Cells.Item(21, 2).Value = 10.4381
$ws.Cells.Item(22, 2).Value = 10.1157
$ws.Cells.Item(22, 5).Value = 16.94730000000002
$ws.Cells.Item(23, 2).Value = 9.196400000000002
$ws.Cells.Item(28, 1).Value = -22.10909999999999
$ws.Cells.Item(29, 1).Value = -21.6798
$ws.Cells.Item(29, 2).Value = 5.367400000000003
$ws.Cells.Item(32, 1).Value = -21.16389999999999
$ws.Cells.Item(34, 2).Value = 9.668600000000005
$ws.Cells.Item(36, 5).Value = 15.92769999999999
$ws.Cells.Item(40, 1).Value = -19.41399999999999
$ws.Cells.Item(42, 2).Value = 10.1164
$ws.Cells.Item(43, 2).Value = 5.685200000000002
$ws.Cells.Item(44, 2).Value = 5.229100000000002
$ws.Cells.Item(45, 2).Value = 4.574400000000003
$ws.Cells.Item(46, 1).Value = -22.01249999999999
$ws.Cells.Item(46, 2).Value = 4.977000000000001
$ws.Cells.Item(46, 5).Value = 16.59709999999999
$ws.Cells.Item(50, 2).Value = 5.144299999999995
$ws.Cells.Item(50, 5).Value = 16.3176
$ws.Cells.Item(51, 1).Value = -22.1633
$ws.Cells.Item(51, 2).Value = 5.7616
$ws.Cells.Item(52, 1).Value = -22.0526
$ws.Cells.Item(57, 1).Value = -22.78180000000001
$ws.Cells.Item(59, 1).Value = -22.1099
$ws.Cells.Item(62, 1).Value = -22.0517
$ws.Cells.Item(66, 1).Value = -21.49639999999999
$ws.Cells.Item(66, 2).Value = 5.023699999999998
$ws.Cells.Item(67, 2).Value = 5.631299999999999
$ws.Cells.Item(73, 1).Value = -19.35139999999998
$ws.Cells.Item(74, 1).Value = -22.12859999999999
$ws.Cells.Item(79, 2).Value = 9.989000000000006
$ws.Cells.Item(84, 2).Value = 5.493199999999999
$ws.Cells.Item(92, 1).Value = -21.41900000000001
$ws.Cells.Item(92, 2).Value = 4.981999999999996
$ws.Cells.Item(95, 5).Value = 17.82440000000001
$ws.Cells.Item(97, 2).Value = 5.339799999999998
$ws.Cells.Item(97, 5).Value = 16.6996
$ws.Cells.Item(100, 1).Value = -22.10479999999999
